# Auto-generated edit script: updates market-price-derived columns (H-N)
# for Leve crafting profit rows across all 8 sheets, per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 364.09525
$ws.Range("J6").Value = 770.1111
$ws.Range("L6").Value = 2310.3333
$ws.Range("N6").Value = -2534.3333

$ws.Range("H8").Value = 342.22223
$ws.Range("J8").Value = 695.25
$ws.Range("L8").Value = 2085.75
$ws.Range("N8").Value = -2363.75

$ws.Range("H29").Value = 3949.4285
$ws.Range("J29").Value = 5437
$ws.Range("L29").Value = 16311
$ws.Range("N29").Value = -16873

$ws.Range("H42").Value = 86.75
$ws.Range("I42").Value = 89
$ws.Range("J42").Value = 80
$ws.Range("K42").Value = 267
$ws.Range("L42").Value = 240
$ws.Range("M42").Value = -37
$ws.Range("N42").Value = -700

$ws.Range("H62").Value = 3756.4285
$ws.Range("J62").Value = 3849.5
$ws.Range("L62").Value = 3849.5
$ws.Range("N62").Value = -5097.5

$ws.Range("H65").Value = 3756.4285
$ws.Range("J65").Value = 3849.5
$ws.Range("L65").Value = 19247.5
$ws.Range("N65").Value = -25487.5

$ws.Range("H86").Value = 2480.0417
$ws.Range("J86").Value = 2419.25
$ws.Range("L86").Value = 2419.25
$ws.Range("N86").Value = -4665.25

$ws.Range("H89").Value = 2480.0417
$ws.Range("J89").Value = 2419.25
$ws.Range("L89").Value = 12096.25
$ws.Range("N89").Value = -23328.25

$ws.Range("H137").Value = 2576.1714
$ws.Range("I137").Value = 1118.3462
$ws.Range("K137").Value = 3355.0386
$ws.Range("M137").Value = -805.0385999999999

$ws.Range("H138").Value = 3970.6848
$ws.Range("I138").Value = 3249
$ws.Range("J138").Value = 4085.238
$ws.Range("K138").Value = 9747
$ws.Range("L138").Value = 12255.714
$ws.Range("M138").Value = -4607
$ws.Range("N138").Value = -22535.714

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17181.799
$ws.Range("I32").Value = 15294.176
$ws.Range("K32").Value = 15294.176
$ws.Range("M32").Value = -15007.176

$ws.Range("H44").Value = 51238.2
$ws.Range("J44").Value = 39047.75
$ws.Range("L44").Value = 39047.75
$ws.Range("N44").Value = -40023.75

$ws.Range("H74").Value = 3099.291
$ws.Range("I74").Value = 2880.0417
$ws.Range("K74").Value = 2880.0417
$ws.Range("M74").Value = -2006.0417

$ws.Range("H77").Value = 3099.291
$ws.Range("I77").Value = 2880.0417
$ws.Range("K77").Value = 14400.2085
$ws.Range("M77").Value = -10032.2085

$ws.Range("H122").Value = 4319.143
$ws.Range("I122").Value = 3267.182
$ws.Range("K122").Value = 9801.545999999998
$ws.Range("M122").Value = -7351.545999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 1640.9333
$ws.Range("I64").Value = 328
$ws.Range("K64").Value = 328
$ws.Range("M64").Value = -103

$ws.Range("H67").Value = 1640.9333
$ws.Range("I67").Value = 328
$ws.Range("K67").Value = 328
$ws.Range("M67").Value = 452

$ws.Range("H80").Value = 396.66666
$ws.Range("I80").Value = 400
$ws.Range("J80").Value = 395
$ws.Range("K80").Value = 400
$ws.Range("L80").Value = 395
$ws.Range("M80").Value = 598
$ws.Range("N80").Value = -2391

$ws.Range("H83").Value = 396.66666
$ws.Range("I83").Value = 400
$ws.Range("J83").Value = 395
$ws.Range("K83").Value = 2000
$ws.Range("L83").Value = 1975
$ws.Range("M83").Value = 2992
$ws.Range("N83").Value = -11959

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 496
$ws.Range("I7").Value = 246.5
$ws.Range("J7").Value = 579.1667
$ws.Range("K7").Value = 246.5
$ws.Range("L7").Value = 579.1667
$ws.Range("M7").Value = -133.5
$ws.Range("N7").Value = -805.1667

$ws.Range("H31").Value = 3663.4033
$ws.Range("I31").Value = 3565.875
$ws.Range("K31").Value = 3565.875
$ws.Range("M31").Value = -3270.875

$ws.Range("H34").Value = 3663.4033
$ws.Range("I34").Value = 3565.875
$ws.Range("K34").Value = 3565.875
$ws.Range("M34").Value = -3363.875

$ws.Range("H132").Value = 14083.077
$ws.Range("I132").Value = 7451.8184
$ws.Range("K132").Value = 22355.4552
$ws.Range("M132").Value = -19825.4552

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 4248.5713
$ws.Range("I14").Value = 4248.5713
$ws.Range("K14").Value = 12745.7139
$ws.Range("M14").Value = -12572.7139

$ws.Range("H112").Value = 4263.5
$ws.Range("I112").Value = 684.6667
$ws.Range("J112").Value = 15000
$ws.Range("K112").Value = 2054.0001
$ws.Range("L112").Value = 45000
$ws.Range("M112").Value = -946.0001000000002
$ws.Range("N112").Value = -47216

$ws.Range("H122").Value = 90879
$ws.Range("I122").Value = 816
$ws.Range("J122").Value = 225973.5
$ws.Range("K122").Value = 7344
$ws.Range("L122").Value = 2033761.5
$ws.Range("M122").Value = -4894
$ws.Range("N122").Value = -2038661.5

$ws.Range("H138").Value = 6836.5835
$ws.Range("I138").Value = 5408.9287
$ws.Range("K138").Value = 16226.7861
$ws.Range("M138").Value = -11086.7861

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 7408.8335
$ws.Range("I3").Value = 8500
$ws.Range("J3").Value = 6863.25
$ws.Range("K3").Value = 8500
$ws.Range("L3").Value = 6863.25
$ws.Range("M3").Value = -8384
$ws.Range("N3").Value = -7095.25

$ws.Range("H39").Value = 50000
$ws.Range("J39").Value = 50000
$ws.Range("L39").Value = 50000
$ws.Range("N39").Value = -51064

$ws.Range("H102").Value = 31717.588
$ws.Range("J102").Value = 12756.615
$ws.Range("L102").Value = 12756.615
$ws.Range("N102").Value = -16000.615

$ws.Range("H105").Value = 62012.273
$ws.Range("J105").Value = 62012.273
$ws.Range("L105").Value = 62012.273
$ws.Range("N105").Value = -69000.273

$ws.Range("H113").Value = 7079.8486
$ws.Range("I113").Value = 9559.764999999999
$ws.Range("K113").Value = 9559.764999999999
$ws.Range("M113").Value = -7389.764999999999

$ws.Range("H120").Value = 0
$ws.Range("J120").Value = 0
$ws.Range("L120").Value = 0
$ws.Range("N120").ClearContents()

$ws.Range("H121").Value = 49998
$ws.Range("J121").Value = 49998
$ws.Range("L121").Value = 49998
$ws.Range("N121").Value = -53492

$ws.Range("H122").Value = 2981.762
$ws.Range("I122").Value = 2890.125
$ws.Range("K122").Value = 8670.375
$ws.Range("M122").Value = -6220.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 933.3333
$ws.Range("J16").Value = 467
$ws.Range("L16").Value = 467
$ws.Range("N16").Value = -807

$ws.Range("H82").Value = 40001616
$ws.Range("I82").Value = 1841.2222
$ws.Range("K82").Value = 1841.2222
$ws.Range("M82").Value = -1480.2222

$ws.Range("H85").Value = 40001616
$ws.Range("I85").Value = 1841.2222
$ws.Range("K85").Value = 1841.2222
$ws.Range("M85").Value = -593.2221999999999

$ws.Range("H122").Value = 4606.9
$ws.Range("I122").Value = 4287.3335
$ws.Range("J122").Value = 5352.5557
$ws.Range("K122").Value = 12862.0005
$ws.Range("L122").Value = 16057.6671
$ws.Range("M122").Value = -10412.0005
$ws.Range("N122").Value = -20957.6671

$ws.Range("H132").Value = 5561.5386
$ws.Range("I132").Value = 3961.5
$ws.Range("K132").Value = 11884.5
$ws.Range("M132").Value = -9354.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 135500
$ws.Range("J94").Value = 135500
$ws.Range("L94").Value = 135500
$ws.Range("N94").Value = -137302

$ws.Range("H103").Value = 18228.2
$ws.Range("J103").Value = 18228.2
$ws.Range("L103").Value = 18228.2
$ws.Range("N103").Value = -20572.2

$ws.Range("H122").Value = 9998.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 9998.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 29995.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -34895.5

$ws.Range("H132").Value = 12416.571
$ws.Range("I132").Value = 9486
$ws.Range("K132").Value = 28458
$ws.Range("M132").Value = -25928

Write-Host "Applied scheduled market-price refresh across all sheets."